# The workbook has a single sheet with a daily price log for
# "Hortaliza, Femacal de La Calera - Espinaca" running from row 2
# through row 399 (row 1 is the header). A new observation needs to be
# inserted as the *first* record of the series (row 297), pushing the
# previously-297..399 records down to 298..400 (the last record ends
# up at row 400). The new record re-uses the surrounding columns of
# the record it displaces (market/region/category/quality/unit/
# origin/etc. are constant across the whole series) but carries its
# own date (Fecha, column D) and volume (Volumen, column J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 297:399 down to 298:400, leaving row 297 blank.
$ws.Rows.Item(297).Insert()

# Seed the new row 297 with the data that is now duplicated in row
# 298 (the record that used to live at row 297 before the shift), then
# overwrite the two cells that actually differ for the new record.
$ws.Range("A298:R298").Copy()
$ws.Range("A297:R297").PasteSpecial()

$ws.Range("D297").Value2 = 44809
$ws.Range("J297").Value2 = 120
